$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (long buildup tickers) for rows 2-8
$ws.Range("B2").Value = "NSE:APOLLO"
$ws.Range("B3").Value = "NSE:BLUESTARCO"
$ws.Range("B4").Value = "NSE:KSL"
$ws.Range("B5").Value = "NSE:NAZARA"
$ws.Range("B6").Value = "NSE:PTCIL"
$ws.Range("B7").Value = "NSE:ROSSELLIND"
$ws.Range("B8").Value = "NSE:RSYSTEMS"

# Update column C (support zone tickers) for rows 2-24
$ws.Range("C2").Value = "NSE:ALBERTDAVD"
$ws.Range("C3").Value = "NSE:APEX"
$ws.Range("C4").Value = "NSE:BFSI"
$ws.Range("C5").Value = "NSE:CORDSCABLE"
$ws.Range("C6").Value = "NSE:DABUR"
$ws.Range("C7").Value = "NSE:DBSTOCKBRO"
$ws.Range("C8").Value = "NSE:ESTER"
$ws.Range("C9").Value = "NSE:EXCELINDUS"
$ws.Range("C10").Value = "NSE:FAZE3Q"
$ws.Range("C11").Value = "NSE:FIBERWEB"
$ws.Range("C12").Value = "NSE:GMMPFAUDLR"
$ws.Range("C13").Value = "NSE:GODREJCP"
$ws.Range("C14").Value = "NSE:GTPL"
$ws.Range("C15").Value = "NSE:GULFPETRO"
$ws.Range("C16").Value = "NSE:HESTERBIO"
$ws.Range("C17").Value = "NSE:HLEGLAS"
$ws.Range("C18").Value = "NSE:JUBLPHARMA"
$ws.Range("C19").Value = "NSE:KABRAEXTRU"
$ws.Range("C20").Value = "NSE:LOVABLE"
$ws.Range("C21").Value = "NSE:ONGC"
$ws.Range("C22").Value = "NSE:PDSL"
$ws.Range("C23").Value = "NSE:PIIND"
$ws.Range("C24").Value = "NSE:SAFARI"

# Update column E (short buildup tickers) for rows 2-6
$ws.Range("E2").Value = "NSE:ADANIENT"
$ws.Range("E3").Value = "NSE:ADANIGREEN"
$ws.Range("E4").Value = "NSE:BDL"
$ws.Range("E5").Value = "NSE:GLENMARK"
$ws.Range("E6").Value = "NSE:LODHA"

# Remove the now-unused trailing rows (25-35), shrinking the sheet's
# used range from A1:F35 down to A1:F24
$ws.Range("A25:F35").EntireRow.Delete()
